$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.679.84"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "1.635.45"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'213.11"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'19.17"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "1.864.91"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.631.98"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "'4.09"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "26.672.48"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "'63.16"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").Value = "'217.76"
$ws.Range("E19").Value = "  +8.05%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'4.30"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").Value = "'9.46"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "'6.21"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D25").Value = "'148.41"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "'6.92"
$ws.Range("E28").Value = "  +5.65%  "
$ws.Range("D29").Value = "'15.44"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'3.31"
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "1.201.00"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("E37").Value = "  +5.98%  "
$ws.Range("D38").Value = "'0.809"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'0.505"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("D43").Value = "'0.794"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "1.771.86"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'92.29"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'1.55"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("D47").Value = "'54.76"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'7.63"
$ws.Range("E49").Value = "  +5.31%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  +0.14%  "
